$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "281.95"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.31%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "28.50"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.79%"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "4.17%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06480"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.56%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.214"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "3.00%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.384"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.58%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.387"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "4.82%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9253"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "5.73%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1533"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.98%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06389"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "24.54%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07539"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.47%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02902"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.39%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.08960"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.00%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001596"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.58%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0006421"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.93%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006105"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.39%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.441"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.99%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.231"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.67%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3186"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.23%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.59%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.066"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "4.16%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.1551"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.13%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04409"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.25%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001193"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.59%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004035"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "3.54%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001257"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "6.66%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001627"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-1.05%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04104"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.57%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006706"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.78%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1221"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-12.31%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002113"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "11.81%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01209"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "3.61%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005653"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "5.73%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.01307"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-29.32%"

Write-Host "Applied all cell updates"